$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Delete the first 4 data rows (spreadsheet rows 2-5), which correspond
    # to the original Cutoff values 0,1,2,3. Remaining rows shift up and
    # keep their B and C values, but column A (Cutoff) is renumbered
    # starting again from 0.
    $ws.Rows("2:5").Delete()

    for ($r = 2; $r -le 16; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 2
    }
}
